$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Obrigatorio) for rows 2 through 10 from "N" to "S"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = "S"
}
